# Change_Log.xlsx edit:
#  - Row 20 on the "Change Log" sheet gets a date (A20) and a real
#    "Changes" note (B20) describing the JAL work, replacing the
#    previously-blank placeholder row.
#  - The active selection moves to D20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")

# A20: date of this log entry (matches the date already used in row 19)
$ws.Range("A20").Value = "10/26/2025"

# B20: changelog text. Build with an explicit LF (Excel's in-cell line
# break) between the "Changes" header and the two bullet lines, with the
# trailing run of spaces preserved exactly as typed by the original author.
$nl = [char]10
$trailingPad = "".PadLeft(231)
$text = "Changes" + $nl + `
        "- MODIFIED: MyMIF.mif, alu.vhd, alu_control.vhd, Controller.vhd, registerfile.vhd" + $nl + `
        "- COMPLETED: Implemented & tested JAL instruction. Added a comment to registerfile.vhd" + `
        $trailingPad

$ws.Range("B20").Value = $text

# Move the active cell/selection to D20
$ws.Range("D20").Select()
